# Commit: "Added handling of common packages."
#
# This edit reflects a regeneration of the class-structure workbook after the
# code generator started treating "common package" types (repository, kafka,
# product, etc.) specially:
#   1. classFields sheet: the fields belonging to StockComponentTests,
#      OrderManageService and StockApp were re-ordered (same data, new order).
#   2. methodNumberOfLines sheet: the default (no-op / framework) constructors
#      are no longer counted, so those five rows are removed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. classFields: reorder the field rows for the affected classes
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("classFields")

$ws3.Cells.Item(6,1).Value  = "pl.piomin.stock.StockComponentTests"
$ws3.Cells.Item(6,2).Value  = "repository"
$ws3.Cells.Item(6,3).Value  = ""
$ws3.Cells.Item(6,4).Value  = "pl.piomin.stock.repository.ProductRepository"

$ws3.Cells.Item(7,1).Value  = "pl.piomin.stock.StockComponentTests"
$ws3.Cells.Item(7,2).Value  = "kafka"
$ws3.Cells.Item(7,3).Value  = "private"
$ws3.Cells.Item(7,4).Value  = "org.springframework.kafka.test.EmbeddedKafkaBroker"

$ws3.Cells.Item(8,1).Value  = "pl.piomin.stock.StockComponentTests"
$ws3.Cells.Item(8,2).Value  = "product"
$ws3.Cells.Item(8,3).Value  = ""
$ws3.Cells.Item(8,4).Value  = "pl.piomin.stock.domain.Product"

$ws3.Cells.Item(9,1).Value  = "pl.piomin.stock.StockComponentTests"
$ws3.Cells.Item(9,2).Value  = "template"
$ws3.Cells.Item(9,3).Value  = "private"
$ws3.Cells.Item(9,4).Value  = "org.springframework.kafka.core.KafkaTemplate"

$ws3.Cells.Item(10,1).Value = "pl.piomin.stock.StockComponentTests"
$ws3.Cells.Item(10,2).Value = "factory"
$ws3.Cells.Item(10,3).Value = "private"
$ws3.Cells.Item(10,4).Value = "org.springframework.kafka.core.ConsumerFactory"

$ws3.Cells.Item(11,1).Value = "pl.piomin.stock.StockComponentTests"
$ws3.Cells.Item(11,2).Value = "LOG"
$ws3.Cells.Item(11,3).Value = "private"
$ws3.Cells.Item(11,4).Value = "org.slf4j.Logger"

$ws3.Cells.Item(12,1).Value = "pl.piomin.stock.service.OrderManageService"
$ws3.Cells.Item(12,2).Value = "SOURCE"
$ws3.Cells.Item(12,3).Value = "private"
$ws3.Cells.Item(12,4).Value = "java.lang.String"

$ws3.Cells.Item(13,1).Value = "pl.piomin.stock.service.OrderManageService"
$ws3.Cells.Item(13,2).Value = "LOG"
$ws3.Cells.Item(13,3).Value = "private"
$ws3.Cells.Item(13,4).Value = "org.slf4j.Logger"

$ws3.Cells.Item(14,1).Value = "pl.piomin.stock.service.OrderManageService"
$ws3.Cells.Item(14,2).Value = "template"
$ws3.Cells.Item(14,3).Value = "private"
$ws3.Cells.Item(14,4).Value = "org.springframework.kafka.core.KafkaTemplate"

$ws3.Cells.Item(15,1).Value = "pl.piomin.stock.service.OrderManageService"
$ws3.Cells.Item(15,2).Value = "repository"
$ws3.Cells.Item(15,3).Value = "private"
$ws3.Cells.Item(15,4).Value = "pl.piomin.stock.repository.ProductRepository"

$ws3.Cells.Item(16,1).Value = "pl.piomin.stock.StockApp"
$ws3.Cells.Item(16,2).Value = "repository"
$ws3.Cells.Item(16,3).Value = "private"
$ws3.Cells.Item(16,4).Value = "pl.piomin.stock.repository.ProductRepository"

$ws3.Cells.Item(17,1).Value = "pl.piomin.stock.StockApp"
$ws3.Cells.Item(17,2).Value = "LOG"
$ws3.Cells.Item(17,3).Value = "private"
$ws3.Cells.Item(17,4).Value = "org.slf4j.Logger"

$ws3.Cells.Item(18,1).Value = "pl.piomin.stock.StockApp"
$ws3.Cells.Item(18,2).Value = "orderManageService"
$ws3.Cells.Item(18,3).Value = ""
$ws3.Cells.Item(18,4).Value = "pl.piomin.stock.service.OrderManageService"

# ---------------------------------------------------------------------------
# 2. methodNumberOfLines: drop the rows for the trivial/common constructors
#    (StockAppTest(), StockComponentTests(), OrderManageService(...),
#    KafkaContainerDevMode(), StockApp()). Delete bottom-up so row numbers
#    of not-yet-deleted rows stay valid.
# ---------------------------------------------------------------------------
$ws11 = $wb.Worksheets.Item("methodNumberOfLines")

$ws11.Rows.Item(22).Delete()
$ws11.Rows.Item(20).Delete()
$ws11.Rows.Item(17).Delete()
$ws11.Rows.Item(13).Delete()
$ws11.Rows.Item(11).Delete()
